{"js": "// Commit: \"patch 'p' in endif\"\n//\n// The template uses a custom Jinja-like tag syntax where block tags are\n// written as \"{%p ... %}\" (e.g. \"{%p if ... %}\", \"{%p for ... %}\",\n// \"{%p endif %}\", \"{%p endfor %}\"). One \"endif\" tag in the retirement\n// account / employer-sponsored section was missing the \"p\" -- it reads\n// \"{% endif %}\" instead of \"{%p endif %}\". This script fixes that one\n// occurrence so it's consistent with every other block tag in the file.\nconst body = context.document.body;\n\n// Search for the literal (non-\"p\") endif tag. There is exactly one such\n// occurrence in the document; every other \"endif\" already reads\n// \"{%p endif %}\" and must be left untouched.\nconst results = body.search(\"{% endif %}\", {\n  matchCase: true,\n  matchWholeWord: false,\n  matchWildcards: false\n});\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"{%p endif %}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Commit: \"patch 'p' in endif\"\n#\n# The template uses a custom Jinja-like tag syntax where block tags are\n# written as \"{%p ... %}\" (e.g. \"{%p if ... %}\", \"{%p for ... %}\",\n# \"{%p endif %}\", \"{%p endfor %}\"). One \"endif\" tag in the retirement\n# account / employer-sponsored section was missing the \"p\" -- it reads\n# \"{% endif %}\" instead of \"{%p endif %}\". This script fixes that one\n# occurrence so it's consistent with every other block tag in the file.\n\n$d = $word.ActiveDocument\n\n# There is exactly one literal \"{% endif %}\" (without the \"p\") in the\n# document; every other \"endif\" already reads \"{%p endif %}\" and must be\n# left untouched, so a plain, case-sensitive literal search is safe and\n# unambiguous.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"{% endif %}\"\n$find.MatchCase = $true\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 0\n\nwhile ($find.Execute()) {\n    $rng.Text = \"{%p endif %}\"\n    $rng.Collapse(0)\n}\n"}
